$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New week column (G) header date - copy the date style from F1 (same numFmt as C1:F1)
$ws.Range("G1").Value = 43939
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# New week data values (4/18/2020 week) for each industry row
$ws.Range("G2").Value = 32490
$ws.Range("G3").Value = 17481
$ws.Range("G4").Value = 224
$ws.Range("G5").Value = 6242
$ws.Range("G6").Value = 18307
$ws.Range("G7").Value = 4625
$ws.Range("G8").Value = 1510
$ws.Range("G9").Value = 24173
$ws.Range("G10").Value = 6825
$ws.Range("G11").Value = 1566
$ws.Range("G12").Value = 10254
$ws.Range("G13").Value = 179
$ws.Range("G14").Value = 11883
$ws.Range("G15").Value = 6673
$ws.Range("G16").Value = 1238
$ws.Range("G17").Value = 2842
$ws.Range("G18").Value = 29782
$ws.Range("G19").Value = 9431
$ws.Range("G20").Value = 7607

# Match the saved selection state from the source file
[void]$ws.Range("M5").Select()
